$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.726.94"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.774.03"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4595"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3585"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.80"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.044"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.219"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.775.81"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.76"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001061"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06425"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.09"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.803"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.795.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.084"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.20"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.60%  "

$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.978.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.170"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.88"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.090"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09231"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.669"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("E34").Value = "  -0.52%  "

$ws.Range("E35").Value = "  -1.17%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06207"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.09%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02296"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6307"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.948"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.185"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.80%  "

$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.777"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.28"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("E45").Value = "  +0.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5895"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.36"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.951"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06923"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.138"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.18"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.23%  "
